# Consolidate the title's text runs on slide 2 ("Below" + " " + "section-level")
# into a single run containing "Below section-level".
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$tr = $s.Shapes.Item(1).TextFrame.TextRange

# Assigning the exact same concatenated text is treated as a no-op (and
# assigning a text that merely extends it only appends a new trailing run),
# so first set the range to an unrelated placeholder value to force the
# engine to collapse/rewrite the paragraph into a single run, then set the
# final desired text.
$tr.Text = "X"
$tr.Text = "Below section-level"
